$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns
# D column values are kept as text (mirrors source inlineStr formatting,
# preserving trailing zeros such as "66.00" or "1.00").

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '46.520.71'
$ws.Range('E2').Value = '  +5.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.292.60'
$ws.Range('E3').Value = '  +2.99%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.09'
$ws.Range('E5').Value = '  +1.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.22'
$ws.Range('E6').Value = '  +11.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.567'
$ws.Range('E7').Value = '  +1.43%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.521'
$ws.Range('E9').Value = '  +5.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.36'
$ws.Range('E10').Value = '  +9.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0787'
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.38'
$ws.Range('E12').Value = '  +5.52%  '
$ws.Range('E13').Value = '  -0.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.645.64'
$ws.Range('E14').Value = '  +3.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.295.97'
$ws.Range('E15').Value = '  +2.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.77'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.811'
$ws.Range('E17').Value = '  +4.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '46.551.34'
$ws.Range('E18').Value = '  +5.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.94'
$ws.Range('E19').Value = '  +5.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0938'
$ws.Range('E20').Value = '  +3.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.99'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '66.00'
$ws.Range('E22').Value = '  +2.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '250.20'
$ws.Range('E23').Value = '  +6.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.89'
$ws.Range('E24').Value = '  +2.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.92'
$ws.Range('E26').Value = '  +4.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '42.15'
$ws.Range('E27').Value = '  +7.75%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.25'
$ws.Range('E28').Value = '  +1.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.85'
$ws.Range('E29').Value = '  +5.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.96'
$ws.Range('E30').Value = '  +3.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.83'
$ws.Range('E31').Value = '  +13.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.60'
$ws.Range('E32').Value = '  +1.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '146.95'
$ws.Range('E33').Value = '  -3.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0792'
$ws.Range('E34').Value = '  +3.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.23'
$ws.Range('E35').Value = '  +13.85%  '
$ws.Range('E36').Value = '  +10.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.117'
$ws.Range('E37').Value = '  +0.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.09'
$ws.Range('E38').Value = '  +19.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.76'
$ws.Range('E39').Value = '  +4.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.97'
$ws.Range('E40').Value = '  +10.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.33'
$ws.Range('E41').Value = '  +5.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0301'
$ws.Range('E42').Value = '  +0.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.97'
$ws.Range('E44').Value = '  +9.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.810.81'
$ws.Range('E45').Value = '  +0.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '87.98'
$ws.Range('E46').Value = '  +20.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.195'
$ws.Range('E47').Value = '  +5.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '72.81'
$ws.Range('E48').Value = '  +7.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.88'
$ws.Range('E49').Value = '  +5.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '95.16'
$ws.Range('E50').Value = '  +0.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.521.07'
$ws.Range('E51').Value = '  +3.03%  '
